# Poster Minor Project - final edits
#  1. Fix missing space: "Abhishek Rathore(0827IT233D01)" -> "Abhishek Rathore (0827IT233D01)"
#  2/3. Merge split "The "/"Comprehensive...System[for Students]"/"...materials."/"...environment."
#       runs into single runs (keep the non-italic formatting of the leading "The " run).
#  4. Rework the "Guided By / Project Coordinator / H.O.D" textbox: bump font size to 31pt,
#     re-punctuate the three lines, and nudge/resize the textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ------------------------------------------------------------------
# 1) Authors line - add a space before the second student's roll number
# ------------------------------------------------------------------
$authors = $s.Shapes.Item(3)
$authorsTr = $authors.TextFrame.TextRange
$para1 = $authorsTr.Paragraphs(1)
$para1.Characters(1, $para1.Text.Length).Text = " Atharv Sharna (0827IT233D02).  Abhishek Rathore (0827IT233D01)."
# This placeholder has spAutoFit - editing its text reflows/resizes the box.
# The source text box keeps its original size, so restore it explicitly.
$authors.Height = 241.7379

# ------------------------------------------------------------------
# 2/3) Abstract textbox - merge the split "The <i>System</i> ..." runs
# ------------------------------------------------------------------
$abstractBox = $s.Shapes.Item(11)
$abstractTr = $abstractBox.TextFrame.TextRange

$paraA = $abstractTr.Paragraphs(4)
$paraA.Characters(1, $paraA.Text.Length).Text = "The Comprehensive Digital Learning System for Students is a web-based platform designed to centralize and simplify access to study materials."

$paraB = $abstractTr.Paragraphs(9)
$paraB.Characters(1, $paraB.Text.Length).Text = "The Comprehensive Digital Learning System transforms traditional learning methods by ensuring that all students have access to quality resources, personalized learning experiences, and an inclusive educational environment."

# ------------------------------------------------------------------
# 4) Acknowledgement textbox - text + font size + position/size tweak
# ------------------------------------------------------------------
$ackBox = $s.Shapes.Item(18)
$ackTr = $ackBox.TextFrame.TextRange

$ackPara1 = $ackTr.Paragraphs(1)
$ackPara1.Characters(1, $ackPara1.Text.Length).Text = "Guided By: Prof Mahendra Verma  Assist. Prof (I.T)"
$ackPara1.Characters(1, $ackPara1.Text.Length).Font.Size = 31

$ackPara2 = $ackTr.Paragraphs(2)
$ackPara2.Characters(1, $ackPara2.Text.Length).Text = "Project Coordinator: Prof. Monika Choudhary, Assist. Prof  (I.T)"
$ackPara2.Characters(1, $ackPara2.Text.Length).Font.Size = 31

$ackPara3 = $ackTr.Paragraphs(3)
# "Prof. Prashant " run - text unchanged, just grow the font
$ackPara3.Characters(1, 15).Font.Size = 31
# "Lakkadwala" run - text unchanged, just grow the font
$ackPara3.Characters(16, 10).Font.Size = 31
# " H.O.D" run - re-punctuate and grow the font
$ackPara3.Characters(26, 6).Text = ", H.O.D (I.T)"
$ackPara3.Characters(26, ", H.O.D (I.T)".Length).Font.Size = 31

# Nudge + resize the acknowledgement textbox (height stays the same).
$ackBox.Left = 1746.0000787401575
$ackBox.Top = 2367.8699212598426
$ackBox.Width = 846.0
